$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.568.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.08%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.447.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.84%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'580.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.32%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'149.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +9.07%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.449.03"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.92%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E9").Value = "'  +0.87%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.97%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +2.67%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.391"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.85%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.036.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.83%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'27.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +7.58%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -0.25%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000175"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.58%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.450.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.55%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'61.702.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.03%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +8.53%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'14.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.82%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.77%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'390.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +3.85%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.567"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.31%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'3.587.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.64%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'72.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.59%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.23%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'5.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.62%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.12%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +3.81%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +3.49%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -11.98%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.66%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'8.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.46%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +1.25%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D36").Value = "'24.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.50%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'7.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.71%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'5.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.52%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +1.18%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'166.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.38%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +3.87%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'27.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +12.45%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +2.06%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +2.35%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.01%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +1.75%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.41%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.599.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +5.77%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -2.91%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'6.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.29%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'23.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.52%  "
$ws.Range("E51").Style = "Normal"
